# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-10 09:19:25
# Applies the attendance recording / reshuffled "Recorded By" list updates
# to the "Session Analysis Results" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: re-order the "Recorded By" (column G) email list for a given row
# ---------------------------------------------------------------------------
function Set-RecordedBy($row, $text) {
    $ws.Range("G$row").Value2 = $text
}

# ---------------------------------------------------------------------------
# Helper: write a literal text value into a cell that would otherwise be
# auto-parsed as a number/percentage by the recalculation engine (e.g.
# "30.1%"), while preserving the cell's original shared style index.
# The leading apostrophe forces text storage (stripped on commit, like
# typing it into Excel); re-pasting formats-only from a same-styled donor
# cell then drops the transient "quote prefix"/number-format style that
# gets created and restores the original style index.
# ---------------------------------------------------------------------------
function Set-TextValue($rangeAddr, $text, $donorAddr) {
    $ws.Range($rangeAddr).Value2 = "'" + $text
    $ws.Range($donorAddr).Copy() | Out-Null
    $ws.Range($rangeAddr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------------
# Row 2 / 24 - reorder recorders list (ANATOMY session 1)
# ---------------------------------------------------------------------------
$g2 = "nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
Set-RecordedBy 2 $g2
Set-RecordedBy 24 $g2

# ---------------------------------------------------------------------------
# Class statistics block (K6:L10)
# ---------------------------------------------------------------------------
$ws.Range("L6").Value2 = 53          # Recorded Sessions
$ws.Range("L7").Value2 = 18          # Missing Sessions
$ws.Range("L8").Value2 = 105         # Pending Sessions
Set-TextValue "L9" "30.1%" "L7"      # Coverage %

# ---------------------------------------------------------------------------
# Row 10 / 32 - reorder recorders list (MICROBIOLOGY)
# ---------------------------------------------------------------------------
$g10 = "esraa.mostafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg"
Set-RecordedBy 10 $g10
Set-RecordedBy 32 $g10

Set-TextValue "L10" "34.7%" "L7"     # Average Attendance %

# ---------------------------------------------------------------------------
# Row 15 / 37 - reorder recorders list (PHARMACOLOGY)
# ---------------------------------------------------------------------------
$g15 = "hana.amr@med.asu.edu.eg, nancy.abdelshafy@med.asu.edu.eg"
Set-RecordedBy 15 $g15
Set-RecordedBy 37 $g15

$ws.Range("O15").Value2 = 8
$ws.Range("P15").Value2 = 1
Set-TextValue "R15" "36.4%" "R16"
Set-TextValue "S15" "31.4%" "S16"

# Row 16 group statistics
$ws.Range("P16").Value2 = 3
$ws.Range("Q16").Value2 = 12

# ---------------------------------------------------------------------------
# Row 18 / 40 - reorder recorders list (PHYSIOLOGY)
# ---------------------------------------------------------------------------
$g18 = "yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
Set-RecordedBy 18 $g18
Set-RecordedBy 40 $g18

# ---------------------------------------------------------------------------
# Row 19 / 41 / 150 / 172 - reorder recorders list (PHYSIOLOGY)
# ---------------------------------------------------------------------------
$g19 = "naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
Set-RecordedBy 19 $g19
Set-RecordedBy 41 $g19
Set-RecordedBy 150 $g19
Set-RecordedBy 172 $g19

# ---------------------------------------------------------------------------
# Row 20 - session got recorded: fill changes from "Not Recorded" (pink)
# to "Recorded" (green). Copy the format from row 2 (a "Recorded" row)
# so the shared style index is reused, then set the new cell values.
# ---------------------------------------------------------------------------
$ws.Range("A2:I2").Copy() | Out-Null
$ws.Range("A20:I20").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("G20").Value2 = "Wafaa.ebida@med.asu.edu.eg"
$ws.Range("H20").Value2 = "23/216"
$ws.Range("I20").Value2 = "Recorded"

$ws.Range("P20").Value2 = 4
$ws.Range("Q20").Value2 = 13

# Row 21 group statistics
$ws.Range("P21").Value2 = 2
$ws.Range("Q21").Value2 = 14

# Row 22 group statistics
$ws.Range("P22").Value2 = 2
$ws.Range("Q22").Value2 = 14

# ---------------------------------------------------------------------------
# Row 25 - reorder recorders list
# ---------------------------------------------------------------------------
Set-RecordedBy 25 "manar.montaser@med.asu.edu.eg, backup@backdoor.com"

# ---------------------------------------------------------------------------
# Row 42 - session fell back to "Not Recorded": fill changes from
# "Pending" (light yellow) to "Not Recorded" (pink). Copy the format from
# row 7 (a "Not Recorded" row) so the shared style index is reused.
# ---------------------------------------------------------------------------
$ws.Range("A7:I7").Copy() | Out-Null
$ws.Range("A42:I42").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I42").Value2 = "Not Recorded"

# ---------------------------------------------------------------------------
# Row 46 / 68 - reorder recorders list
# ---------------------------------------------------------------------------
$g46 = "hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg"
Set-RecordedBy 46 $g46
Set-RecordedBy 68 $g46

# ---------------------------------------------------------------------------
# Row 54 / 76 / 98 / 120 - reorder recorders list
# ---------------------------------------------------------------------------
$g54 = "Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, merna.said@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, maimustafa@med.asu.edu.eg"
Set-RecordedBy 54 $g54
Set-RecordedBy 76 $g54
Set-RecordedBy 98 $g54
Set-RecordedBy 120 $g54

# ---------------------------------------------------------------------------
# Row 58 / 80 - reorder recorders list
# ---------------------------------------------------------------------------
$g58 = "Amr-Saeed@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
Set-RecordedBy 58 $g58
Set-RecordedBy 80 $g58

# ---------------------------------------------------------------------------
# Row 62 / 84 - reorder recorders list
# ---------------------------------------------------------------------------
$g62 = "aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
Set-RecordedBy 62 $g62
Set-RecordedBy 84 $g62

# ---------------------------------------------------------------------------
# Row 63 / 85 - reorder recorders list
# ---------------------------------------------------------------------------
$g63 = "Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
Set-RecordedBy 63 $g63
Set-RecordedBy 85 $g63

# ---------------------------------------------------------------------------
# Row 64 / 86 - reorder recorders list
# ---------------------------------------------------------------------------
$g64 = "youstina.magdy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
Set-RecordedBy 64 $g64
Set-RecordedBy 86 $g64

# ---------------------------------------------------------------------------
# Row 81 - reorder recorders list
# ---------------------------------------------------------------------------
Set-RecordedBy 81 "enas.omran@med.asu.edu.eg, user@user.com, Walaa.h.ghanima@med.asu.edu.eg"

# ---------------------------------------------------------------------------
# Row 90 / 112 - reorder recorders list
# ---------------------------------------------------------------------------
$g90 = "manar.montaser@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
Set-RecordedBy 90 $g90
Set-RecordedBy 112 $g90

# ---------------------------------------------------------------------------
# Row 96 / 118 - reorder recorders list
# ---------------------------------------------------------------------------
$g96 = "mariam.noureldin@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
Set-RecordedBy 96 $g96
Set-RecordedBy 118 $g96

# ---------------------------------------------------------------------------
# Row 106 / 128 - reorder recorders list
# ---------------------------------------------------------------------------
$g106 = "youstina.magdy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
Set-RecordedBy 106 $g106
Set-RecordedBy 128 $g106

# ---------------------------------------------------------------------------
# Row 107 / 129 - reorder recorders list
# ---------------------------------------------------------------------------
$g107 = "neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
Set-RecordedBy 107 $g107
Set-RecordedBy 129 $g107

# ---------------------------------------------------------------------------
# Row 125 - "Pending" -> "Not Recorded" (copy format from row 7)
# ---------------------------------------------------------------------------
$ws.Range("A7:I7").Copy() | Out-Null
$ws.Range("A125:I125").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("I125").Value2 = "Not Recorded"

# ---------------------------------------------------------------------------
# Row 134 - reorder recorders list
# ---------------------------------------------------------------------------
Set-RecordedBy 134 "asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# ---------------------------------------------------------------------------
# Row 135 - "Pending" -> "Not Recorded" (copy format from row 7)
# ---------------------------------------------------------------------------
$ws.Range("A7:I7").Copy() | Out-Null
$ws.Range("A135:I135").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("I135").Value2 = "Not Recorded"

# ---------------------------------------------------------------------------
# Row 142 / 164 - reorder recorders list
# ---------------------------------------------------------------------------
$g142 = "esraa.mostafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, merna.said@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg"
Set-RecordedBy 142 $g142
Set-RecordedBy 164 $g142

# ---------------------------------------------------------------------------
# Row 151 / 173 - reorder recorders list
# ---------------------------------------------------------------------------
$g151 = "yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
Set-RecordedBy 151 $g151
Set-RecordedBy 173 $g151

# ---------------------------------------------------------------------------
# Row 156 - reorder recorders list
# ---------------------------------------------------------------------------
Set-RecordedBy 156 "alshimaa.atef@med.asu.edu.egm, majorelle.magdy@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"

# ---------------------------------------------------------------------------
# Row 157 - "Pending" -> "Not Recorded" (copy format from row 7)
# ---------------------------------------------------------------------------
$ws.Range("A7:I7").Copy() | Out-Null
$ws.Range("A157:I157").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("I157").Value2 = "Not Recorded"

Write-Host "Attendance report sync complete."
